$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, shifting the existing data (rows 20-142) down
# by one row (rows become 21-143).
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new daily price record.
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C20").Value = "Arica y Parinacota"
$ws.Range("D20").Value = 44749
$ws.Range("E20").Value = 15
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100108
$ws.Range("H20").Value = "Tropicales y subtropicales"
$ws.Range("I20").Value = 100108002
$ws.Range("J20").Value = "Mango"
$ws.Range("K20").Value = "Sin especificar"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = 9000
$ws.Range("O20").Value = 10000
$ws.Range("P20").Value = 9500
$ws.Range("Q20").Value = "$/bandeja 4 kilos"
$ws.Range("R20").Value = "Brasil"
$ws.Range("S20").Value = 2375
$ws.Range("T20").Value = 4
